$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update area_lid value
$ws.Range("B5").Value = 1

# Update dQ_cool formula to use named ranges
$ws.Range("B8").Formula = "=k_lid*area_lid*(set_temp-dt_cool)/thickness_lid"

# Update c_water value
$ws.Range("B10").Value = 1

# Update fluid_mass value
$ws.Range("B11").Value = 1

# Update Q_cool formula to use named ranges
$ws.Range("B12").Formula = "=fluid_mass*c_water*dt_cool"

# Update the active selection to B8
$ws.Range("B8").Select()
